$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 9 (shifts "Andel heltidsstilling" and everything
# below it down by one row) and populate it with the new "deltidsstilling"
# (part-time position) variable.
$ws.Rows.Item(9).Insert()

$ws.Cells.Item(9, 1).Value = "Andel deltidsstilling"
$ws.Cells.Item(9, 2).Value = "deltidsstilling"
$ws.Cells.Item(9, 3).Value = "snitt_as_num_single"
$ws.Cells.Item(9, 4).Value = "Andel deltid"

# The table ("Table1") covered the inserted row's range but does not
# auto-grow on a plain row insert, so resize it explicitly to include
# the new row (A1:E30 -> A1:E31).
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E31"))

$ws.Range("C10").Select()
